# Add support for timevals (#83)
# Replace the plain 2018-2021 year values in the "Data" sheet (column C)
# with quarter-style time values (2000Q1..2000Q4), and move the active
# sheet/tab selection from "Codelists" to "Data".

$wb = $excel.ActiveWorkbook

# --- Update the "Data" worksheet (4th sheet / rId4) ---------------------
$wsData = $wb.Worksheets.Item(4)

$quarters = @("2000Q1", "2000Q2", "2000Q3", "2000Q4")
for ($r = 2; $r -le 85; $r++) {
    $q = $quarters[($r - 2) % 4]
    $wsData.Cells.Item($r, 3).Value = $q
}

# --- Move the tab/selection state ---------------------------------------
# Previously "Codelists" (3rd sheet) was the selected/active tab with
# selection C11; now "Data" (4th sheet) becomes active with selection E16,
# and "Codelists" loses its tabSelected flag (selection stays C11).
$wsCodelists = $wb.Worksheets.Item(3)
$wsCodelists.Activate() | Out-Null
$wsCodelists.Range("C11").Select() | Out-Null

$wsData.Activate() | Out-Null
$wsData.Range("E16").Select() | Out-Null
